# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns
# to match the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.958.11"
$ws.Cells.Item(2, 5).Value = "  +0.25%  "

$ws.Cells.Item(3, 4).Value = "1.877.45"
$ws.Cells.Item(3, 5).Value = "  -0.89%  "

$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$ws.Cells.Item(5, 4).Value = "'0.7406"
$ws.Cells.Item(5, 5).Value = "  -4.31%  "

$ws.Cells.Item(6, 4).Value = "'242.23"
$ws.Cells.Item(6, 5).Value = "  -0.85%  "

$ws.Cells.Item(7, 5).Value = "  -0.04%  "

$ws.Cells.Item(8, 4).Value = "'0.3154"
$ws.Cells.Item(8, 5).Value = "  +0.82%  "

$ws.Cells.Item(9, 4).Value = "'0.07173"
$ws.Cells.Item(9, 5).Value = "  -0.68%  "

$ws.Cells.Item(10, 4).Value = "'24.80"
$ws.Cells.Item(10, 5).Value = "  -3.49%  "

$ws.Cells.Item(11, 4).Value = "'0.08430"
$ws.Cells.Item(11, 5).Value = "  -4.79%  "

$ws.Cells.Item(12, 4).Value = "'0.7551"
$ws.Cells.Item(12, 5).Value = "  -2.12%  "

$ws.Cells.Item(13, 4).Value = "'5.420"
$ws.Cells.Item(13, 5).Value = "  -0.20%  "

$ws.Cells.Item(14, 4).Value = "1.872.22"
$ws.Cells.Item(14, 5).Value = "  -3.80%  "

$ws.Cells.Item(15, 4).Value = "'92.85"
$ws.Cells.Item(15, 5).Value = "  -1.70%  "

$ws.Cells.Item(16, 4).Value = "29.959.36"
$ws.Cells.Item(16, 5).Value = "  -0.14%  "

$ws.Cells.Item(17, 4).Value = "'6.103"
$ws.Cells.Item(17, 5).Value = "  -1.48%  "

$ws.Cells.Item(18, 5).Value = "  -2.14%  "

$ws.Cells.Item(19, 4).Value = "'243.73"
$ws.Cells.Item(19, 5).Value = "  -0.73%  "

$ws.Cells.Item(20, 4).Value = "'0.000007836"
$ws.Cells.Item(20, 5).Value = "  -0.33%  "

$ws.Cells.Item(21, 4).Value = "'0.9998"
$ws.Cells.Item(21, 5).Value = "  -0.01%  "

$ws.Cells.Item(22, 4).Value = "2.116.12"
$ws.Cells.Item(22, 5).Value = "  -3.58%  "

$ws.Cells.Item(23, 4).Value = "'7.994"
$ws.Cells.Item(23, 5).Value = "  -2.03%  "

$ws.Cells.Item(24, 4).Value = "'1.001"
$ws.Cells.Item(24, 5).Value = "  +0.02%  "

$ws.Cells.Item(25, 5).Value = "  -2.30%  "

$ws.Cells.Item(26, 4).Value = "'9.325"

$ws.Cells.Item(27, 4).Value = "'164.75"
$ws.Cells.Item(27, 5).Value = "  +1.43%  "

$ws.Cells.Item(28, 4).Value = "'18.68"

$ws.Cells.Item(29, 5).Value = "  +0.00%  "

$ws.Cells.Item(30, 4).Value = "'1.479"
$ws.Cells.Item(30, 5).Value = "  +3.55%  "

$ws.Cells.Item(31, 4).Value = "'4.613"
$ws.Cells.Item(31, 5).Value = "  +1.25%  "

$ws.Cells.Item(32, 4).Value = "'1.533"
$ws.Cells.Item(32, 5).Value = "  -0.60%  "

$ws.Cells.Item(33, 4).Value = "'4.297"
$ws.Cells.Item(33, 5).Value = "  +4.45%  "

$ws.Cells.Item(34, 4).Value = "'0.05341"
$ws.Cells.Item(34, 5).Value = "  -2.72%  "

$ws.Cells.Item(35, 5).Value = "  -0.62%  "

$ws.Cells.Item(36, 4).Value = "'0.7598"
$ws.Cells.Item(36, 5).Value = "  +1.08%  "

$ws.Cells.Item(37, 4).Value = "'1.000"
$ws.Cells.Item(37, 5).Value = "  +0.03%  "

$ws.Cells.Item(38, 4).Value = "'2.701"
$ws.Cells.Item(38, 5).Value = "  -0.58%  "

$ws.Cells.Item(39, 5).Value = "  +0.04%  "

$ws.Cells.Item(40, 4).Value = "'2.753"
$ws.Cells.Item(40, 5).Value = "  -1.24%  "

$ws.Cells.Item(41, 4).Value = "'0.4493"
$ws.Cells.Item(41, 5).Value = "  -0.31%  "

$ws.Cells.Item(42, 4).Value = "1.113.77"
$ws.Cells.Item(42, 5).Value = "  +1.98%  "

$ws.Cells.Item(43, 4).Value = "'6.149"
$ws.Cells.Item(43, 5).Value = "  +1.93%  "

$ws.Cells.Item(44, 4).Value = "'72.68"
$ws.Cells.Item(44, 5).Value = "  -1.58%  "

$ws.Cells.Item(45, 4).Value = "'0.8642"

$ws.Cells.Item(46, 5).Value = "  +0.14%  "

$ws.Cells.Item(47, 4).Value = "'103.27"
$ws.Cells.Item(47, 5).Value = "  +0.48%  "

$ws.Cells.Item(48, 4).Value = "'7.707"
$ws.Cells.Item(48, 5).Value = "  +1.24%  "

$ws.Cells.Item(49, 4).Value = "'1.849"
$ws.Cells.Item(49, 5).Value = "  -2.08%  "

$ws.Cells.Item(50, 5).Value = "  +4.14%  "

$ws.Cells.Item(51, 4).Value = "2.013.35"
$ws.Cells.Item(51, 5).Value = "  -2.63%  "
